$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.35323633438204638
$ws.Range("A2").Value = -0.0099999993387811514
$ws.Range("A3").Value = -0.0089999993313032434
$ws.Range("A4").Value = -0.011999999828987029
$ws.Range("A5").Value = -0.0059999993446693267
$ws.Range("A6").Value = -0.0059999993282850994
$ws.Range("A7").Value = -0.019999999202505947
$ws.Range("A8").Value = -0.019999999201340657
$ws.Range("A9").Value = -0.0059999993268755603
$ws.Range("A10").Value = -0.005999999328402339
$ws.Range("A11").Value = -0.0044999993419203577
$ws.Range("A12").Value = -0.0059999993289556741
$ws.Range("A13").Value = -0.0059999993347181757
$ws.Range("A14").Value = -0.01199999928374762
$ws.Range("A15").Value = 0.065696628804108315
$ws.Range("A16").Value = -0.0059999993390782436
$ws.Range("A17").Value = -0.0059999993365673632
$ws.Range("A18").Value = -0.0089999993092053643
$ws.Range("A19").Value = -0.066219943749434851
$ws.Range("A20").Value = -0.0089999993402010148
$ws.Range("A21").Value = -0.0089999993393847788
$ws.Range("A22").Value = -0.0089999993388110155
$ws.Range("A23").Value = -0.0089999993254332722
$ws.Range("A24").Value = -0.041999999020404744
$ws.Range("A25").Value = -0.04199999901530127
$ws.Range("A26").Value = -0.0059999993261747875
$ws.Range("A27").Value = -0.0059999993233224025
$ws.Range("A28").Value = -0.005999999311223192
$ws.Range("A29").Value = -0.011999999250168258
$ws.Range("A30").Value = -0.019999999175076333
$ws.Range("A31").Value = -0.0013244178816496088
$ws.Range("A32").Value = -0.020999999161005256
$ws.Range("A33").Value = -0.0059999992946737635
